$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 525
$ws.Range("F6").Value = 1124
$ws.Range("F11").Value = 1120
$ws.Range("F13").Value = 92
$ws.Range("F14").Value = 763
$ws.Range("F15").Value = 797
$ws.Range("F17").Value = 35
$ws.Range("F18").Value = 61
$ws.Range("F19").Value = 665
$ws.Range("F21").Value = 1707
$ws.Range("F22").Value = 2158
$ws.Range("F23").Value = 582
$ws.Range("F25").Value = 1835
$ws.Range("F26").Value = 284
$ws.Range("F27").Value = 2669
$ws.Range("F30").Value = 668
$ws.Range("F31").Value = 128
$ws.Range("F34").Value = 923
$ws.Range("F35").Value = 1643
$ws.Range("F38").Value = 528
$ws.Range("F39").Value = 136
$ws.Range("F40").Value = 110
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 9
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 525
$ws.Range("F7").Value = 1124
$ws.Range("F12").Value = 1120
$ws.Range("F14").Value = 763
$ws.Range("F15").Value = 797
$ws.Range("F20").Value = 35
$ws.Range("F22").Value = 61
$ws.Range("F23").Value = 665
$ws.Range("F25").Value = 1707
$ws.Range("F26").Value = 2158
$ws.Range("F27").Value = 582
$ws.Range("F29").Value = 9
$ws.Range("F31").Value = 2669
$ws.Range("F38").Value = 668
$ws.Range("F39").Value = 128
$ws.Range("F42").Value = 923
$ws.Range("F43").Value = 1643
$ws.Range("F46").Value = 528
$ws.Range("F47").Value = 136
$ws.Range("F48").Value = 110
